# Working version of the tool
#
# Update the propellant/tank parameter values and move the view back to
# where the user was last working (selection + scroll position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fuel/oxidizer mass (kg): 180.8 -> 174
$ws.Range("B4").Value = 174

# Fibres volume fraction (-): 0.25 -> 0.275
$ws.Range("B15").Value = 0.275

# Scroll so row 10 is the top visible row, then move the selection to E17
# (previously the window was scrolled to the top and G14 was selected).
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E17").Select() | Out-Null
